$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 numeric columns with new TPM-derived values
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.779790999999999
$ws.Range("N2").Value = 14.339373
$ws.Range("Q2").Value = 0.111675036924
$ws.Range("R2").Value = 1.005075332316
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove the now-obsolete second data row (previously "Resolving-Mac")
$ws.Rows.Item(3).Delete()
